$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Paragraph 1: "Does CORDIC work without modification in 4th
#    quadrant as well ? "  ->  append two spaces then bold "Yes"
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$ins1 = $d.Range($r1.Start, $r1.End - 1)
$ins1.Collapse(0)
$ins1.InsertAfter("  ")
$ins1.Collapse(0)
$boldStart1 = $ins1.End
$boldText1 = "Yes"
$ins1.InsertAfter($boldText1)
$boldEnd1 = $boldStart1 + $boldText1.Length
$boldRange1 = $d.Range($boldStart1, $boldEnd1)
$boldRange1.Bold = 1

# ---------------------------------------------------------------------
# 2) Paragraph 2: "Look at theta = 90 degrees (2's complement needed ??)"
#    -> append space then bold
#    "Not needed but there is a wrong value before the right value is displayed"
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$ins2 = $d.Range($r2.Start, $r2.End - 1)
$ins2.Collapse(0)
$ins2.InsertAfter(" ")
$ins2.Collapse(0)
$boldStart2 = $ins2.End
$boldText2 = "Not needed but there is a wrong value before the right value is displayed"
$ins2.InsertAfter($boldText2)
$boldEnd2 = $boldStart2 + $boldText2.Length
$boldRange2 = $d.Range($boldStart2, $boldEnd2)
$boldRange2.Bold = 1

# ---------------------------------------------------------------------
# 3) The empty paragraph right after "Have a separate module
#    get_conv_target_angle(...)" becomes a new, indented, bold
#    paragraph: "All four quadrants working"
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5.LeftIndent = 72   # 72 points == 1440 twips
$r5 = $p5.Range
$ins5 = $d.Range($r5.Start, $r5.Start)
$newText5 = "All four quadrants working"
$ins5.InsertAfter($newText5)
$p5.Range.Font.Bold = 1
